$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# New averaged-intensity run added three "Spiral" sampling schemes.
# The table grows from 16 to 19 rows: row 10 now reports the
# (recomputed) Gaussian-Quadrature numbers, rows 11-13 are the three
# new Spiral schemes, and the remaining schemes (NoRotation-tilt60deg
# through HexGrid-60degTilt5degRes) shift down into rows 14-19.
# ------------------------------------------------------------------

# Row 10: Gaussian-Quadrature (recomputed values)
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 1.18355117195812
$ws.Range("D10").Value = 0.5960752793557353
$ws.Range("E10").Value = 1.047915988402299
$ws.Range("F10").Value = 1.18355117195812
$ws.Range("G10").Value = 0.7943672782176859
$ws.Range("H10").Value = 1.124909215325272
$ws.Range("I10").Value = 1.092202312283082
$ws.Range("J10").Value = 0.5960752793557353
$ws.Range("K10").Value = 0.8219956338790171
$ws.Range("L10").Value = 1.002773402918569
$ws.Range("M10").Value = 0.9731702075903658

# Row 11: Spiral-90deg-10rot-5space (new)
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 0.8715835538134338
$ws.Range("D11").Value = 1.014577400179688
$ws.Range("E11").Value = 1.072895145672294
$ws.Range("F11").Value = 0.8715835538134338
$ws.Range("G11").Value = 0.9114504995053088
$ws.Range("H11").Value = 1.310849725526502
$ws.Range("I11").Value = 1.012997855574666
$ws.Range("J11").Value = 1.014577400179688
$ws.Range("K11").Value = 1.043736272925991
$ws.Range("L11").Value = 0.9576599133697123
$ws.Range("M11").Value = 1.032392363378649

# Row 12: Spiral-90deg-15rot-5space (new)
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 0.8720107334434821
$ws.Range("D12").Value = 1.015579854310746
$ws.Range("E12").Value = 1.072455887706028
$ws.Range("F12").Value = 0.8720107334434821
$ws.Range("G12").Value = 0.911950648142667
$ws.Range("H12").Value = 1.309806607583377
$ws.Range("I12").Value = 1.012571661214339
$ws.Range("J12").Value = 1.015579854310746
$ws.Range("K12").Value = 1.044017871008387
$ws.Range("L12").Value = 0.9580143022259349
$ws.Range("M12").Value = 1.03239589873344

# Row 13: Spiral-90deg-10rot-3space (new)
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 0.8720130697531544
$ws.Range("D13").Value = 1.014271392900869
$ws.Range("E13").Value = 1.072886620186914
$ws.Range("F13").Value = 0.8720130697531544
$ws.Range("G13").Value = 0.9114446638162848
$ws.Range("H13").Value = 1.310140465724916
$ws.Range("I13").Value = 1.012846226983344
$ws.Range("J13").Value = 1.014271392900869
$ws.Range("K13").Value = 1.043579006543891
$ws.Range("L13").Value = 0.957796038148523
$ws.Range("M13").Value = 1.03226707322758

# Row 14: NoRotation-tilt60deg (was row 10)
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 0.8051360000000007
$ws.Range("D14").Value = 0.5146320000000005
$ws.Range("E14").Value = 1.283668000000001
$ws.Range("F14").Value = 0.8051360000000007
$ws.Range("G14").Value = 0.5394719999999995
$ws.Range("H14").Value = 2.114767999999995
$ws.Range("I14").Value = 1.192755999999999
$ws.Range("J14").Value = 0.5146320000000005
$ws.Range("K14").Value = 0.899150000000001
$ws.Range("L14").Value = 0.8521430000000008
$ws.Range("M14").Value = 1.075071999999999

# Row 15: Rotation-NoTilt (was row 11)
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 0.65
$ws.Range("D15").Value = 0.1
$ws.Range("E15").Value = 1.5
$ws.Range("F15").Value = 0.65
$ws.Range("G15").Value = 0.2198625
$ws.Range("H15").Value = 2.910137500000002
$ws.Range("I15").Value = 1.35
$ws.Range("J15").Value = 0.1
$ws.Range("K15").Value = 0.8
$ws.Range("L15").Value = 0.725
$ws.Range("M15").Value = 1.121666666666667

# Row 16: Rotation-60detTilt (was row 12)
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 0.8045082411008023
$ws.Range("D16").Value = 0.4685234790399995
$ws.Range("E16").Value = 1.289028794368001
$ws.Range("F16").Value = 0.8045082411008023
$ws.Range("G16").Value = 0.5440728799232009
$ws.Range("H16").Value = 2.090021382348802
$ws.Range("I16").Value = 1.197373526630398
$ws.Range("J16").Value = 0.4685234790399995
$ws.Range("K16").Value = 0.8787761367040003
$ws.Range("L16").Value = 0.8416421889024013
$ws.Range("M16").Value = 1.065588050568534

# Row 17: HexGrid-90degTilt5degRes (was row 13) -- new row, copy formatting from row 16 first
$ws.Range("A16").Copy() | Out-Null
$ws.Range("A17").PasteSpecial(-4122) | Out-Null
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.9885873240038706
$ws.Range("D17").Value = 0.9965694279424336
$ws.Range("E17").Value = 0.9969546953521777
$ws.Range("F17").Value = 0.9885873240038706
$ws.Range("G17").Value = 0.9918684008992457
$ws.Range("H17").Value = 0.9998573886316411
$ws.Range("I17").Value = 0.994928375453311
$ws.Range("J17").Value = 0.9965694279424336
$ws.Range("K17").Value = 0.9967620616473056
$ws.Range("L17").Value = 0.9926746928255882
$ws.Range("M17").Value = 0.9947942687137799

# Row 18: HexGrid-90degTilt22p5degRes (was row 14) -- new row, copy formatting from row 16 first
$ws.Range("A16").Copy() | Out-Null
$ws.Range("A18").PasteSpecial(-4122) | Out-Null
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 1.007944356450144
$ws.Range("D18").Value = 1.050266089062521
$ws.Range("E18").Value = 0.9684491221893732
$ws.Range("F18").Value = 1.007944356450144
$ws.Range("G18").Value = 1.01854047218831
$ws.Range("H18").Value = 0.9534279518733793
$ws.Range("I18").Value = 0.9824754194620022
$ws.Range("J18").Value = 1.050266089062521
$ws.Range("K18").Value = 1.009357605625947
$ws.Range("L18").Value = 1.008650981038045
$ws.Range("M18").Value = 0.9968505685376217

# Row 19: HexGrid-60degTilt5degRes (was row 15) -- new row, copy formatting from row 16 first
$ws.Range("A16").Copy() | Out-Null
$ws.Range("A19").PasteSpecial(-4122) | Out-Null
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 0.9749403871273977
$ws.Range("D19").Value = 1.153434291898012
$ws.Range("E19").Value = 0.9507920037967438
$ws.Range("F19").Value = 0.9749403871273977
$ws.Range("G19").Value = 1.08440653595496
$ws.Range("H19").Value = 0.8627221319899996
$ws.Range("I19").Value = 0.9516158174577607
$ws.Range("J19").Value = 1.153434291898012
$ws.Range("K19").Value = 1.052113147847378
$ws.Range("L19").Value = 1.013526767487388
$ws.Range("M19").Value = 0.9963185280374788
